# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Kraken_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the source diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 54.083332
$ws.Range("I11").Value = 54.083332
$ws.Range("K11").Value = 54.083332
$ws.Range("M11").Value = 85.916668
$ws.Range("H38").Value = 72.14286
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()  # was -1746
$ws.Range("H42").Value = 39
$ws.Range("I42").Value = 41.25
$ws.Range("J42").Value = 34.5
$ws.Range("K42").Value = 123.75
$ws.Range("L42").Value = 103.5
$ws.Range("M42").Value = 106.25
$ws.Range("N42").Value = -563.5
$ws.Range("H43").Value = 4428
$ws.Range("I43").Value = 5666
$ws.Range("K43").Value = 5666
$ws.Range("M43").Value = -5597
$ws.Range("H58").Value = 2963
$ws.Range("J58").Value = 4000
$ws.Range("L58").Value = 12000
$ws.Range("N58").Value = -12300
$ws.Range("H70").Value = 17883.143
$ws.Range("I70").Value = 3867.25
$ws.Range("K70").Value = 11601.75
$ws.Range("M70").Value = -11331.75
$ws.Range("H73").Value = 17883.143
$ws.Range("I73").Value = 3867.25
$ws.Range("K73").Value = 11601.75
$ws.Range("M73").Value = -10665.75
$ws.Range("H98").Value = 2586.125
$ws.Range("I98").Value = 2500.8333
$ws.Range("J98").Value = 2842
$ws.Range("K98").Value = 2500.8333
$ws.Range("L98").Value = 2842
$ws.Range("M98").Value = -1002.8333
$ws.Range("N98").Value = -5838
$ws.Range("H122").Value = 2586.125
$ws.Range("I122").Value = 2500.8333
$ws.Range("J122").Value = 2842
$ws.Range("K122").Value = 7502.499899999999
$ws.Range("L122").Value = 8526
$ws.Range("M122").Value = -5052.499899999999
$ws.Range("N122").Value = -13426

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 771.1111
$ws.Range("I97").Value = 780
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 780
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -284
$ws.Range("N97").Value = -1692
$ws.Range("H102").Value = 1865.6666
$ws.Range("I102").Value = 1865.6666
$ws.Range("K102").Value = 1865.6666
$ws.Range("M102").Value = -243.6666
$ws.Range("H110").Value = 980.6667
$ws.Range("I110").Value = 980.6667
$ws.Range("K110").Value = 980.6667
$ws.Range("M110").Value = 1064.3333
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 23439.445
$ws.Range("I86").Value = 1789.6
$ws.Range("J86").Value = 50501.75
$ws.Range("K86").Value = 1789.6
$ws.Range("L86").Value = 50501.75
$ws.Range("M86").Value = -666.5999999999999
$ws.Range("N86").Value = -52747.75
$ws.Range("H89").Value = 23439.445
$ws.Range("I89").Value = 1789.6
$ws.Range("J89").Value = 50501.75
$ws.Range("K89").Value = 8948
$ws.Range("L89").Value = 252508.75
$ws.Range("M89").Value = -3332
$ws.Range("N89").Value = -263740.75
$ws.Range("H94").Value = 2275.087
$ws.Range("I94").Value = 1515.85
$ws.Range("K94").Value = 1515.85
$ws.Range("M94").Value = -1064.85
$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 4000
$ws.Range("K99").Value = 4000
$ws.Range("M99").Value = -2502
$ws.Range("H105").Value = 3699.8
$ws.Range("I105").Value = 3499.75
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 3499.75
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -1752.75
$ws.Range("N105").Value = -7994
$ws.Range("H107").Value = 2650.1428
$ws.Range("I107").Value = 1410.2
$ws.Range("K107").Value = 1410.2
$ws.Range("M107").Value = 509.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51316
$ws.Range("H122").Value = 1594.5714
$ws.Range("I122").Value = 1534.6
$ws.Range("J122").Value = 1744.5
$ws.Range("K122").Value = 4603.799999999999
$ws.Range("L122").Value = 5233.5
$ws.Range("M122").Value = -2153.799999999999
$ws.Range("N122").Value = -10133.5
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 999.5
$ws.Range("I3").Value = 999.5
$ws.Range("K3").Value = 2998.5
$ws.Range("M3").Value = -2886.5
$ws.Range("H26").Value = 130.33333
$ws.Range("I26").Value = 45.5
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 136.5
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = 151.5
$ws.Range("N26").Value = -1476
$ws.Range("H68").Value = 1002.4
$ws.Range("J68").Value = 1078
$ws.Range("L68").Value = 3234
$ws.Range("N68").Value = -4856
$ws.Range("H71").Value = 1002.4
$ws.Range("J71").Value = 1078
$ws.Range("L71").Value = 9702
$ws.Range("N71").Value = -17814
$ws.Range("H133").Value = 4533
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 4533
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 13599
$ws.Range("M133").ClearContents()  # was -19436.0005
$ws.Range("N133").Value = -23719
$ws.Range("H134").Value = 450
$ws.Range("I134").Value = 450
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1350
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 3720
$ws.Range("N134").ClearContents()  # was -11637
$ws.Range("H136").Value = 1343.5
$ws.Range("I136").Value = 1343.5
$ws.Range("K136").Value = 4030.5
$ws.Range("M136").Value = 1069.5
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -19200
$ws.Range("H138").Value = 2454.1428
$ws.Range("I138").Value = 2363.1667
$ws.Range("K138").Value = 7089.500100000001
$ws.Range("M138").Value = -1949.500100000001
$ws.Range("H141").Value = 9833
$ws.Range("I141").Value = 1999.6666
$ws.Range("J141").Value = 33333
$ws.Range("K141").Value = 5998.9998
$ws.Range("L141").Value = 99999
$ws.Range("M141").Value = -818.9997999999996
$ws.Range("N141").Value = -110359

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1257.6666
$ws.Range("I97").Value = 1314.5294
$ws.Range("J97").Value = 291
$ws.Range("K97").Value = 1314.5294
$ws.Range("L97").Value = 291
$ws.Range("M97").Value = -818.5293999999999
$ws.Range("N97").Value = -1283
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("H122").Value = 1384.5
$ws.Range("I122").Value = 1361.4
$ws.Range("K122").Value = 4084.2
$ws.Range("M122").Value = -1634.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 252
$ws.Range("I9").Value = 252
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 252
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -28
$ws.Range("N9").ClearContents()  # was -1448
$ws.Range("H40").Value = 10187.5625
$ws.Range("I40").Value = 10366.733
$ws.Range("K40").Value = 10366.733
$ws.Range("M40").Value = -10230.733
$ws.Range("H93").Value = 7950
$ws.Range("I93").Value = 9359.200000000001
$ws.Range("J93").Value = 904
$ws.Range("K93").Value = 9359.200000000001
$ws.Range("L93").Value = 904
$ws.Range("M93").Value = -8111.200000000001
$ws.Range("N93").Value = -3400
$ws.Range("H132").Value = 10750
$ws.Range("I132").Value = 10750
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 32250
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -29720
$ws.Range("N132").ClearContents()  # was -35060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 392999
$ws.Range("J118").Value = 392999
$ws.Range("L118").Value = 392999
$ws.Range("N118").Value = -396313
$ws.Range("H141").Value = 99994
$ws.Range("J141").Value = 99994
$ws.Range("L141").Value = 99994
$ws.Range("N141").Value = -110354
